$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.10622251033783
$ws.Range("B1").Value = 1.45796012878418
$ws.Range("C1").Value = 9.084759712219238
$ws.Range("D1").Value = 2.388652801513672
$ws.Range("E1").Value = 1.281911015510559
